# Add a new "localdb" command-category to the #system sheet:
#  - shift columns N..AC one column to the right (-> O..AD) to make room for
#    a brand-new "N" column holding the localdb function list
#  - shift the "target" category list (column A, rows 14..29) down one row
#    to make room (alphabetically) for the new "localdb" entry at A14
#  - populate the new column N (rows 1..7) with the localdb header + its
#    6 functions
#  - fix up every defined name whose range moved, and add the new "localdb"
#    defined name

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------
# 1) Shift the "target" list in column A down by one row, rows 14-29 -> 15-30,
#    so we can insert "localdb" (alphabetically, between "json" and "macro")
#    at A14. Done manually (bottom-up) so only column A is touched -- the
#    built-in Range.Insert shifts the whole row's cells, which we don't want.
# ---------------------------------------------------------------------
for ($r = 29; $r -ge 14; $r--) {
    $v = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r + 1, 1).Value2 = $v
}
$ws.Cells.Item(14, 1).Value2 = "localdb"

# ---------------------------------------------------------------------
# 2) Shift columns N:AC one column to the right (-> O:AD) across every row,
#    freeing up column N for the new localdb function list.
# ---------------------------------------------------------------------
$ws.Columns("N:N").Insert()

# ---------------------------------------------------------------------
# 3) Populate the new column N with the localdb header + functions.
# ---------------------------------------------------------------------
$ws.Cells.Item(1, 14).Value2 = "localdb"
$ws.Cells.Item(2, 14).Value2 = "cloneTable(var,source,target)"
$ws.Cells.Item(3, 14).Value2 = "dropTables(var,tables)"
$ws.Cells.Item(4, 14).Value2 = "exportCSV(sql,output)"
$ws.Cells.Item(5, 14).Value2 = "importRecords(var,sourceDb,sql,table)"
$ws.Cells.Item(6, 14).Value2 = "purge(var)"
$ws.Cells.Item(7, 14).Value2 = "runSQLs(var,sqls)"

# ---------------------------------------------------------------------
# 4) Fix up defined names whose target range moved because of the column
#    shift, plus "target" (row count grew) and the brand-new "localdb" name.
# ---------------------------------------------------------------------
$wb.Names.Item("mail").RefersTo       = "='#system'!`$P`$2:`$P`$2"
$wb.Names.Item("number").RefersTo     = "='#system'!`$Q`$2:`$Q`$16"
$wb.Names.Item("pdf").RefersTo        = "='#system'!`$R`$2:`$R`$16"
$wb.Names.Item("rdbms").RefersTo      = "='#system'!`$S`$2:`$S`$7"
$wb.Names.Item("redis").RefersTo      = "='#system'!`$T`$2:`$T`$10"
$wb.Names.Item("sms").RefersTo        = "='#system'!`$U`$2:`$U`$2"
$wb.Names.Item("sound").RefersTo      = "='#system'!`$V`$2:`$V`$5"
$wb.Names.Item("ssh").RefersTo        = "='#system'!`$W`$2:`$W`$9"
$wb.Names.Item("step").RefersTo       = "='#system'!`$X`$2:`$X`$4"
$wb.Names.Item("target").RefersTo     = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("web").RefersTo        = "='#system'!`$Y`$2:`$Y`$127"
$wb.Names.Item("webalert").RefersTo   = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo  = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo         = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo   = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo        = "='#system'!`$AD`$2:`$AD`$21"
$wb.Names.Item("macro").RefersTo      = "='#system'!`$O`$2:`$O`$4"

$wb.Names.Add("localdb", "='#system'!`$N`$2:`$N`$7")
